$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) for rows 2 through 15 from 45186 to 45188,
# preserving the existing cell style/format.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
